$wb = $excel.ActiveWorkbook

# ============================================================
# "Sheet2" -> "Modem Test": insert a Timestamp column at the
# front, add a new data row, and a couple of formatted-but-
# empty trailing rows (mirrors the author's raw edit).
# ============================================================
$modem = $wb.Worksheets.Item("Sheet2")
$modem.Name = "Modem Test"

# remember the hyperlinked cells' text before the column shifts,
# so we can restore it exactly after repairing the hyperlinks
$hlOrigVals = @{}
$hlOrigVals["C2"] = $modem.Range("C2").Value2
$hlOrigVals["C3"] = $modem.Range("C3").Value2
$hlOrigVals["C4"] = $modem.Range("C4").Value2
$hlOrigVals["C5"] = $modem.Range("C5").Value2
$hlOrigVals["C6"] = $modem.Range("C6").Value2
$hlOrigVals["C7"] = $modem.Range("C7").Value2
$hlOrigVals["C8"] = $modem.Range("C8").Value2

# shift A:I -> B:J, making room for the new Timestamp column
$modem.Columns.Item(1).Insert()

# the hyperlink anchors don't follow the column insert automatically,
# so rebuild them at their new column (D) location
$modem.Hyperlinks.Delete()
$hlTargets = @(
  @{cell="D2"; orig="C2"; url="https://github.com/Qrist0ph/pirmicboard_david/blob/ccd93bc5e59aa4cae0e1f3a7a5dfe34171ab0504/unittests/7080gconnect/tmp_main_works_on_gpstracker.py"},
  @{cell="D3"; orig="C3"; url="https://github.com/Qrist0ph/pirmicboard_david/blob/ccd93bc5e59aa4cae0e1f3a7a5dfe34171ab0504/unittests/7080gconnect/tmp_main_works_on_gpstracker.py"},
  @{cell="D4"; orig="C4"; url="https://github.com/Qrist0ph/pirmicboard_david/blob/ccd93bc5e59aa4cae0e1f3a7a5dfe34171ab0504/unittests/7080gconnect/tmp_main_works_on_gpstracker.py"},
  @{cell="D5"; orig="C5"; url="https://github.com/Qrist0ph/pirmicboard_david/blob/a2ef3bf7a16e352b20008353330c1dd8c7383065/unittests/7080gconnect/main.py"},
  @{cell="D6"; orig="C6"; url="https://github.com/Qrist0ph/pirmicboard_david/blob/a2ef3bf7a16e352b20008353330c1dd8c7383065/unittests/7080gconnect/main.py"},
  @{cell="D7"; orig="C7"; url="https://github.com/Qrist0ph/pirmicboard_david/blob/5e8b1c0b49c3b408311a3a79990f6b82077685f9/unittests/7080gconnect/main.py"},
  @{cell="D8"; orig="C8"; url="https://github.com/Qrist0ph/pirmicboard_david/blob/5e8b1c0b49c3b408311a3a79990f6b82077685f9/unittests/7080gconnect/main.py"}
)
foreach ($t in $hlTargets) {
  $modem.Hyperlinks.Add($modem.Range($t.cell), $t.url, "", "", $t.url)
  $modem.Range($t.cell).Value2 = $hlOrigVals[$t.orig]
  $modem.Range($t.cell).Style = "Hyperlink"
}

# new Timestamp column + the new "90 sek bug" test row
$modem.Range("A1").Value = "Timestamp"
$modem.Range("A9").Value = 45196.387499999997
$modem.Range("A9").NumberFormat = "m/d/yy h:mm"

$modem.Range("B9").Value = "Board getauscht"
$modem.Range("C9").Value = 2
$modem.Range("D9").Value = "main.py"
$modem.Range("H9").Value = "vebunden nach ca 30 run"

# leftover formatted-but-empty cells from the original edit
$modem.Range("B10").NumberFormat = "h:mm AM/PM"
$modem.Range("B11").NumberFormat = "m/d/yy h:mm"

# restore approximate column widths after the insert shifted everything right
$modem.Columns.Item(1).ColumnWidth = 12.83
$modem.Columns.Item(2).ColumnWidth = 40.83
$modem.Columns.Item(5).ColumnWidth = 20.5
$modem.Columns.Item(6).ColumnWidth = 19.39
$modem.Columns.Item(8).ColumnWidth = 45.83
$modem.Columns.Item(9).ColumnWidth = 78.94

# ============================================================
# New "PIR Test" sheet (right after "Modem Test")
# ============================================================
$pir = $wb.Worksheets.Add($null, $modem)
$pir.Name = "PIR Test"

$pir.Range("A1").Value = "Timestamp"
$pir.Range("B1").Value = "Board"
$pir.Range("C1").Value = "Fesnel"
$pir.Range("D1").Value = "Firmware"
$pir.Range("E1").Value = "90 Sek Bug"
$pir.Range("F1").Value = "Szenario"
$pir.Range("G1").Value = "Result"

$pir.Range("A2").Value = 45196.385416666664
$pir.Range("A2").NumberFormat = "m/d/yy h:mm"
$pir.Range("B2").Value = 3
$pir.Range("C2").Value = "SR 501"
$pir.Range("G2").Value = "PIR scheint sehr empfindlich"

$pir.Range("A3").Value = 45196.388888888891
$pir.Range("A3").NumberFormat = "m/d/yy h:mm"
$pir.Range("B3").Value = 2
$pir.Range("C3").Value = "SR 501"
$pir.Range("E3").Value = "yes"
$pir.Range("G3").Value = "PIR sieht gut aus"

$pir.Range("A4").Value = 45196.45
$pir.Range("A4").NumberFormat = "m/d/yy h:mm"
$pir.Range("B4").Value = 2
$pir.Range("E4").Value = "NEIN, wenn PIR unplugged"
$pir.Range("F4").Value = "PIR unplugged"

$pir.Range("E5").Select()

$pir.Columns.Item(1).ColumnWidth = 13.83
$pir.Columns.Item(5).ColumnWidth = 18.61
$pir.Columns.Item(6).ColumnWidth = 8.72
$pir.Columns.Item(7).ColumnWidth = 23.17

# make "PIR Test" the active sheet/tab (target workbook has activeTab=2)
$pir.Activate()
